$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2921").Value = "2022-01-05"; $ws.Range("B2921").Value = 4.1; $ws.Range("C2921").Value = "30:14"
$ws.Range("A2922").Value = "2022-01-06"; $ws.Range("B2922").Value = 6; $ws.Range("C2922").Value = "42:55"
$ws.Range("A2923").Value = "2022-01-10"; $ws.Range("B2923").Value = 4.05; $ws.Range("C2923").Value = "28:40"
$ws.Range("A2924").Value = "2022-01-12"; $ws.Range("B2924").Value = 10; $ws.Range("C2924").Value = "1:08:35"
$ws.Range("A2925").Value = "2022-01-19"; $ws.Range("B2925").Value = 4.03; $ws.Range("C2925").Value = "29:10"
$ws.Range("A2926").Value = "2022-01-22"; $ws.Range("B2926").Value = 3.16; $ws.Range("C2926").Value = "23:53"
$ws.Range("A2927").Value = "2022-01-25"; $ws.Range("B2927").Value = 4.07; $ws.Range("C2927").Value = "27:16"
$ws.Range("A2928").Value = "2022-02-25"; $ws.Range("B2928").Value = 3.45; $ws.Range("C2928").Value = "29:13"
$ws.Range("A2929").Value = "2022-03-06"; $ws.Range("B2929").Value = 4.38; $ws.Range("C2929").Value = "32:58"
$ws.Range("A2930").Value = "2022-03-10"; $ws.Range("B2930").Value = 4.2; $ws.Range("C2930").Value = "34:37"
$ws.Range("A2931").Value = "2022-03-14"; $ws.Range("B2931").Value = 4.44; $ws.Range("C2931").Value = "35:21"
$ws.Range("A2932").Value = "2022-03-17"; $ws.Range("B2932").Value = 5.55; $ws.Range("C2932").Value = "41:06"
$ws.Range("A2933").Value = "2022-03-22"; $ws.Range("B2933").Value = 2.75; $ws.Range("C2933").Value = "23:17"
$ws.Range("A2934").Value = "2022-03-26"; $ws.Range("B2934").Value = 5.03; $ws.Range("C2934").Value = "40:27"
$ws.Range("A2935").Value = "2022-04-04"; $ws.Range("B2935").Value = 3.11; $ws.Range("C2935").Value = "23:25"
$ws.Range("A2936").Value = "2022-04-06"; $ws.Range("B2936").Value = 4.4; $ws.Range("C2936").Value = "30:26"
$ws.Range("A2937").Value = "2022-04-09"; $ws.Range("B2937").Value = 5.29; $ws.Range("C2937").Value = "36:20"
$ws.Range("A2938").Value = "2022-04-11"; $ws.Range("B2938").Value = 4.1; $ws.Range("C2938").Value = "28:15"
$ws.Range("A2939").Value = "2022-04-16"; $ws.Range("B2939").Value = 6.16; $ws.Range("C2939").Value = "41:37"
$ws.Range("A2940").Value = "2022-04-22"; $ws.Range("B2940").Value = 4.01; $ws.Range("C2940").Value = "29:47"
$ws.Range("A2941").Value = "2022-04-25"; $ws.Range("B2941").Value = 5; $ws.Range("C2941").Value = "34:35"
$ws.Range("A2942").Value = "2022-04-28"; $ws.Range("B2942").Value = 6.07; $ws.Range("C2942").Value = "42:29"
$ws.Range("A2943").Value = "2022-05-01"; $ws.Range("B2943").Value = 7.3; $ws.Range("C2943").Value = "55:43"
$ws.Range("A2944").Value = "2022-05-05"; $ws.Range("B2944").Value = 6; $ws.Range("C2944").Value = "42:39"
$ws.Range("A2945").Value = "2022-05-08"; $ws.Range("B2945").Value = 5.55; $ws.Range("C2945").Value = "39:47"
$ws.Range("A2946").Value = "2022-05-11"; $ws.Range("B2946").Value = 4.18; $ws.Range("C2946").Value = "29:47"
$ws.Range("A2947").Value = "2022-05-18"; $ws.Range("B2947").Value = 6.11; $ws.Range("C2947").Value = "43:03"
$ws.Range("A2948").Value = "2022-05-24"; $ws.Range("B2948").Value = 4.1; $ws.Range("C2948").Value = "29:19"
$ws.Range("A2949").Value = "2022-05-26"; $ws.Range("B2949").Value = 4.93; $ws.Range("C2949").Value = "39:29"
$ws.Range("A2950").Value = "2022-06-03"; $ws.Range("B2950").Value = 4.89; $ws.Range("C2950").Value = "35:08"
$ws.Range("A2951").Value = "2022-06-07"; $ws.Range("B2951").Value = 4.05; $ws.Range("C2951").Value = "29:01"
$ws.Range("A2952").Value = "2022-06-14"; $ws.Range("B2952").Value = 4.07; $ws.Range("C2952").Value = "29:43"
$ws.Range("A2953").Value = "2022-06-16"; $ws.Range("B2953").Value = 4.34; $ws.Range("C2953").Value = "31:48"
$ws.Range("A2954").Value = "2022-07-02"; $ws.Range("B2954").Value = 3.11; $ws.Range("C2954").Value = "24:18"
$ws.Range("A2955").Value = "2022-07-04"; $ws.Range("B2955").Value = 4.19; $ws.Range("C2955").Value = "32:37"
$ws.Range("A2956").Value = "2022-07-14"; $ws.Range("B2956").Value = 4.03; $ws.Range("C2956").Value = "28:55"
$ws.Range("A2957").Value = "2022-07-22"; $ws.Range("B2957").Value = 5.03; $ws.Range("C2957").Value = "36:28"
$ws.Range("A2958").Value = "2022-07-25"; $ws.Range("B2958").Value = 3.7; $ws.Range("C2958").Value = "28:03"
$ws.Range("A2959").Value = "2022-07-29"; $ws.Range("B2959").Value = 5.29; $ws.Range("C2959").Value = "39:23"
$ws.Range("A2960").Value = "2022-07-31"; $ws.Range("B2960").Value = 3.15; $ws.Range("C2960").Value = "23:10"
$ws.Range("A2961").Value = "2022-08-03"; $ws.Range("B2961").Value = 3.11; $ws.Range("C2961").Value = "19:49"
$ws.Range("A2962").Value = "2022-08-03"; $ws.Range("B2962").Value = 2.04
$ws.Range("A2963").Value = "2022-08-07"; $ws.Range("B2963").Value = 5.13; $ws.Range("C2963").Value = "41:35"
$ws.Range("A2964").Value = "2022-08-10"; $ws.Range("B2964").Value = 6.11; $ws.Range("C2964").Value = "48:27"
$ws.Range("A2965").Value = "2022-08-12"; $ws.Range("B2965").Value = 5.86; $ws.Range("C2965").Value = "49:25"
$ws.Range("A2966").Value = "2022-08-14"; $ws.Range("B2966").Value = 5.22; $ws.Range("C2966").Value = "39:24"
$ws.Range("A2967").Value = "2022-08-17"; $ws.Range("B2967").Value = 6.04; $ws.Range("C2967").Value = "44:17"
$ws.Range("A2968").Value = "2022-08-19"; $ws.Range("B2968").Value = 5.18; $ws.Range("C2968").Value = "38:23"
$ws.Range("A2969").Value = "2022-08-22"; $ws.Range("B2969").Value = 5.02; $ws.Range("C2969").Value = "35:44"
$ws.Range("A2970").Value = "2022-08-24"; $ws.Range("B2970").Value = 4.15; $ws.Range("C2970").Value = "32:52"
$ws.Range("A2971").Value = "2022-08-26"; $ws.Range("B2971").Value = 6.22; $ws.Range("C2971").Value = "45:06"
$ws.Range("A2972").Value = "2022-08-29"; $ws.Range("B2972").Value = 5.88; $ws.Range("C2972").Value = "43:22"
$ws.Range("A2973").Value = "2022-08-31"; $ws.Range("B2973").Value = 7.21; $ws.Range("C2973").Value = "56:29"
$ws.Range("A2974").Value = "2022-09-02"; $ws.Range("B2974").Value = 3.2; $ws.Range("C2974").Value = "24:27"
$ws.Range("A2975").Value = "2022-09-04"; $ws.Range("B2975").Value = 7.59; $ws.Range("C2975").Value = "1:08:13"
$ws.Range("A2976").Value = "2022-09-06"; $ws.Range("B2976").Value = 4.28; $ws.Range("C2976").Value = "30:41"
$ws.Range("A2977").Value = "2022-09-07"; $ws.Range("B2977").Value = 6.15; $ws.Range("C2977").Value = "48:29"
$ws.Range("A2978").Value = "2022-09-09"; $ws.Range("B2978").Value = 6.43; $ws.Range("C2978").Value = "48:13"
$ws.Range("A2979").Value = "2022-09-11"; $ws.Range("B2979").Value = 8.19; $ws.Range("C2979").Value = "1:03:49"
$ws.Range("A2980").Value = "2022-09-12"; $ws.Range("B2980").Value = 5.02; $ws.Range("C2980").Value = "34:42"
$ws.Range("A2981").Value = "2022-09-14"; $ws.Range("B2981").Value = 9.15; $ws.Range("C2981").Value = "1:06:59"
$ws.Range("A2982").Value = "2022-09-16"; $ws.Range("B2982").Value = 6.65; $ws.Range("C2982").Value = "49:19"
$ws.Range("A2983").Value = "2022-09-17"; $ws.Range("B2983").Value = 7.04; $ws.Range("C2983").Value = "49:21"
$ws.Range("A2984").Value = "2022-09-19"; $ws.Range("B2984").Value = 6.61; $ws.Range("C2984").Value = "48:59"
$ws.Range("A2985").Value = "2022-09-21"; $ws.Range("B2985").Value = 7.8; $ws.Range("C2985").Value = "58:13"
$ws.Range("A2986").Value = "2022-09-23"; $ws.Range("B2986").Value = 5.66; $ws.Range("C2986").Value = "38:08"
$ws.Range("A2987").Value = "2022-09-24"; $ws.Range("B2987").Value = 10.21; $ws.Range("C2987").Value = "1:16:13"
$ws.Range("A2988").Value = "2022-09-26"; $ws.Range("B2988").Value = 4.3; $ws.Range("C2988").Value = "31:34"
$ws.Range("A2989").Value = "2022-09-28"; $ws.Range("B2989").Value = 8.93; $ws.Range("C2989").Value = "1:05:22"
$ws.Range("A2990").Value = "2022-09-30"; $ws.Range("B2990").Value = 4.42; $ws.Range("C2990").Value = "31:55"
$ws.Range("A2991").Value = "2022-10-01"; $ws.Range("B2991").Value = 5.71; $ws.Range("C2991").Value = "42:18"
$ws.Range("A2992").Value = "2022-10-03"; $ws.Range("B2992").Value = 7.35; $ws.Range("C2992").Value = "52:23"
$ws.Range("A2993").Value = "2022-10-05"; $ws.Range("B2993").Value = 9.26; $ws.Range("C2993").Value = "1:07:41"
$ws.Range("A2994").Value = "2022-10-07"; $ws.Range("B2994").Value = 4.02; $ws.Range("C2994").Value = "28:38"
$ws.Range("A2995").Value = "2022-10-09"; $ws.Range("B2995").Value = 6.55; $ws.Range("C2995").Value = "42:06"
$ws.Range("A2996").Value = "2022-10-10"; $ws.Range("B2996").Value = 7; $ws.Range("C2996").Value = "50:42"
$ws.Range("A2997").Value = "2022-10-12"; $ws.Range("B2997").Value = 11.11; $ws.Range("C2997").Value = "1:23:24"
$ws.Range("A2998").Value = "2022-10-14"; $ws.Range("B2998").Value = 7.3; $ws.Range("C2998").Value = "47:41"
$ws.Range("A2999").Value = "2022-10-15"; $ws.Range("B2999").Value = 4.24; $ws.Range("C2999").Value = "29:16"
$ws.Range("A3000").Value = "2022-10-17"; $ws.Range("B3000").Value = 6; $ws.Range("C3000").Value = "43:30"
$ws.Range("A3001").Value = "2022-10-19"; $ws.Range("B3001").Value = 7.7; $ws.Range("C3001").Value = "54:30"
$ws.Range("A3002").Value = "2022-10-20"; $ws.Range("B3002").Value = 12.08; $ws.Range("C3002").Value = "1:24:39"
$ws.Range("A3003").Value = "2022-10-23"; $ws.Range("B3003").Value = 6.23; $ws.Range("C3003").Value = "46:39"
$ws.Range("A3004").Value = "2022-10-24"; $ws.Range("B3004").Value = 8.78; $ws.Range("C3004").Value = "1:02:05"
$ws.Range("A3005").Value = "2022-10-26"; $ws.Range("B3005").Value = 9.88; $ws.Range("C3005").Value = "1:11:54"
$ws.Range("A3006").Value = "2022-10-29"; $ws.Range("B3006").Value = 10.02; $ws.Range("C3006").Value = "1:22:14"
$ws.Range("A3007").Value = "2022-10-31"; $ws.Range("B3007").Value = 6.07; $ws.Range("C3007").Value = "46:09"
$ws.Range("A3008").Value = "2022-11-02"; $ws.Range("B3008").Value = 8.04; $ws.Range("C3008").Value = "55:46"
$ws.Range("A3009").Value = "2022-11-04"; $ws.Range("B3009").Value = 8.23; $ws.Range("C3009").Value = "54:28"
$ws.Range("A3010").Value = "2022-11-06"; $ws.Range("B3010").Value = 11.67; $ws.Range("C3010").Value = "1:18:11"
$ws.Range("A3011").Value = "2022-11-08"; $ws.Range("B3011").Value = 8.18; $ws.Range("C3011").Value = "58:45"
$ws.Range("A3012").Value = "2022-11-09"; $ws.Range("B3012").Value = 8.35; $ws.Range("C3012").Value = "1:01:02"
$ws.Range("A3013").Value = "2022-11-11"; $ws.Range("B3013").Value = 9.11; $ws.Range("C3013").Value = "1:03:35"
$ws.Range("A3014").Value = "2022-11-13"; $ws.Range("B3014").Value = 12.11; $ws.Range("C3014").Value = "1:23:55"
$ws.Range("A3015").Value = "2022-11-16"; $ws.Range("B3015").Value = 5.01; $ws.Range("C3015").Value = "38:57"
$ws.Range("A3016").Value = "2022-11-18"; $ws.Range("B3016").Value = 6.29; $ws.Range("C3016").Value = "43:58"
$ws.Range("A3017").Value = "2022-11-19"; $ws.Range("B3017").Value = 3.11; $ws.Range("C3017").Value = "17:08"
$ws.Range("A3018").Value = "2022-11-19"; $ws.Range("B3018").Value = 4.33
$ws.Range("A3019").Value = "2022-11-21"; $ws.Range("B3019").Value = 6.5; $ws.Range("C3019").Value = "50:18"
$ws.Range("A3020").Value = "2022-11-23"; $ws.Range("B3020").Value = 6.91; $ws.Range("C3020").Value = "1:01:59"
$ws.Range("A3021").Value = "2022-11-24"; $ws.Range("B3021").Value = 3.32; $ws.Range("C3021").Value = "17:58"
$ws.Range("A3022").Value = "2022-11-24"; $ws.Range("B3022").Value = 4.81
$ws.Range("A3023").Value = "2022-11-26"; $ws.Range("B3023").Value = 6.22; $ws.Range("C3023").Value = "40:49"
$ws.Range("A3024").Value = "2022-11-29"; $ws.Range("B3024").Value = 5; $ws.Range("C3024").Value = "36:23"

$ws.Range("B3025").Select()

